# Introduction to Spring Framework Testing
# Adds two new "Spring Web Reactive" slides (mock request/response support,
# then WebTestClient) right after the current slide 7, pushing the
# pre-existing trailing blank slides further down the deck.

$p = $ppt.ActivePresentation

function Set-RunFont {
    param($paragraphRange, [string]$needle, [string]$fontName)
    $txt = $paragraphRange.Text
    $idx = $txt.IndexOf($needle)
    if ($idx -lt 0) {
        throw "needle not found: [$needle] in [$txt]"
    }
    $sub = $paragraphRange.Characters($idx + 1, $needle.Length)
    $sub.Font.Name = $fontName
}

# Same as Set-RunFont, but also formats the single space that follows
# the needle as its own run (matches how PowerPoint split the package
# names from the following space when the font change was applied).
function Set-RunFontPlusSpace {
    param($paragraphRange, [string]$needle, [string]$fontName)
    $txt = $paragraphRange.Text
    $idx = $txt.IndexOf($needle)
    if ($idx -lt 0) {
        throw "needle not found: [$needle] in [$txt]"
    }
    $sub = $paragraphRange.Characters($idx + 1, $needle.Length)
    $sub.Font.Name = $fontName
    $spaceRun = $paragraphRange.Characters($idx + 1 + $needle.Length, 1)
    $spaceRun.Font.Name = $fontName
}

function Add-Paragraph {
    param($textRange, [string]$text)
    [void]$textRange.InsertAfter("`r" + $text)
}

# ---------------------------------------------------------------------
# The deck currently ends with three blank "Title and Content" slides
# (slides 7, 8 and 9). Duplicating slide 7 four times in a row places
# the new slides immediately after it, which reproduces the target
# p:sldId order 256,257,258,259,260,261,262,265,266,267,268,263,264.
# ---------------------------------------------------------------------
$src = $p.Slides.Item(7)
[void]$src.Duplicate()
[void]$p.Slides.Item(8).Duplicate()
[void]$p.Slides.Item(9).Duplicate()
[void]$p.Slides.Item(10).Duplicate()

# ---------------------------------------------------------------------
# Slide 7 (was blank): "Spring Web Reactive" - mock request/response
# ---------------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
$slide7.Shapes.Item(1).TextFrame.TextRange.Text = "Spring Web Reactive"

$body7 = $slide7.Shapes.Item(2).TextFrame.TextRange
$slide7.Shapes.Item(2).TextFrame.AutoSize = 2

$body7.Text = "The mock implementations of ServerHttpRequest and ServerHttpResponse are provided for testing WebFlux applications are provided in the org.springframework.mock.http.server.reactive package. "
Add-Paragraph $body7 "The org.springframework.mock.web.server package contains a mock ServerWebExchange that depends on those mock implementations."
Add-Paragraph $body7 "Both MockServerHttpRequest and MockServerHttpResponse extend from the same abstract base classes as server-specific implementations and share behavior with them. "
Add-Paragraph $body7 "For example, a mock request is immutable once created, but you can use the mutate() method from ServerHttpRequest to create a modified instance."

$body7.Paragraphs(2, 1).IndentLevel = 2
$body7.Paragraphs(4, 1).IndentLevel = 2

$para1 = $body7.Paragraphs(1, 1)
Set-RunFont $para1 "ServerHttpRequest" "Courier"
Set-RunFont $para1 "ServerHttpResponse" "Courier"
Set-RunFontPlusSpace $para1 "org.springframework.mock.http.server.reactive" "Courier"

$para2 = $body7.Paragraphs(2, 1)
Set-RunFontPlusSpace $para2 "org.springframework.mock.web.server" "Courier"
Set-RunFont $para2 "ServerWebExchange" "Courier"

$para3 = $body7.Paragraphs(3, 1)
Set-RunFont $para3 "MockServerHttpRequest" "Courier"
Set-RunFont $para3 "MockServerHttpResponse" "Courier"

$para4 = $body7.Paragraphs(4, 1)
Set-RunFont $para4 "mutate()" "Courier"
Set-RunFont $para4 "ServerHttpRequest" "Courier"

# ---------------------------------------------------------------------
# Slide 8 (was blank): "Spring Web Reactive" - WebTestClient
# ---------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$slide8.Shapes.Item(1).TextFrame.TextRange.Text = "Spring Web Reactive"

$body8 = $slide8.Shapes.Item(2).TextFrame.TextRange
$body8.Text = "In order for the mock response to properly implement the write contract and return a write completion handle (that is, Mono<Void>), it by default uses a Flux with cache().then(), which buffers the data and makes it available for assertions in tests. "
Add-Paragraph $body8 "Applications can set a custom write function (for example, to test an infinite stream)."
Add-Paragraph $body8 "The WebTestClient builds on the mock request and response to provide support for testing WebFlux applications without an HTTP server. "
Add-Paragraph $body8 "The client can also be used for end-to-end tests with a running server."
Add-Paragraph $body8 ""
Add-Paragraph $body8 ""

$body8.Paragraphs(2, 1).IndentLevel = 2
$body8.Paragraphs(4, 1).IndentLevel = 2

$para8_3 = $body8.Paragraphs(3, 1)
Set-RunFont $para8_3 "WebTestClient" "Courier"

Write-Output ("Slides: " + $p.Slides.Count)
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    $t = $sl.Shapes.Item(1).TextFrame.TextRange.Text
    Write-Output ("idx=$i sldId=" + $sl.SlideID + " title=[" + $t + "]")
}
